$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.569.42"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "1.878.94"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("D5").Value = "'313.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").Value = "'0.4798"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "'0.3786"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("D9").Value = "'0.07396"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").Value = "'0.9428"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "'20.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.88%  "
$ws.Range("D12").Value = "'0.07876"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.40%  "
$ws.Range("D13").Value = "1.861.62"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "'5.456"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").Value = "'6.616"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.42%  "
$ws.Range("D16").Value = "'91.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "'0.000008982"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "'15.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.06%  "
$ws.Range("D21").Value = "27.598.52"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").Value = "'5.151"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").Value = "'10.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").Value = "'1.967"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").Value = "'153.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").Value = "'18.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").Value = "'2.031"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("D28").Value = "'116.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "'5.008"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").Value = "'0.08936"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "'3.323"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").Value = "'1.217"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("D33").Value = "'4.628"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").Value = "'0.7521"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").Value = "'2.704"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.06%  "
$ws.Range("D36").Value = "'0.02079"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.66%  "
$ws.Range("D37").Value = "'1.124"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("D38").Value = "'0.05314"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "'3.011"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("D40").Value = "'0.5383"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.40%  "
$ws.Range("D41").Value = "'7.116"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").Value = "'0.1524"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").Value = "'8.470"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.22%  "
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("D45").Value = "'0.4849"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.19%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "'1.666"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.02%  "
$ws.Range("D48").Value = "'103.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D49").Value = "'67.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("D50").Value = "'0.06113"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").Value = "'0.9022"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.85%  "
